$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-26 18:18:13"
$ws.Range("O2").Value = "5.9 °C"
$ws.Range("E3").Value = "2026-02-26 18:18:16"
$ws.Range("E4").Value = "2026-02-26 18:18:18"
$ws.Range("E5").Value = "2026-02-26 18:18:21"
$ws.Range("E6").Value = "2026-02-26 18:18:23"
$ws.Range("E7").Value = "2026-02-26 18:18:25"
$ws.Range("E8").Value = "2026-02-26 18:18:28"
$ws.Range("E9").Value = "2026-02-26 18:18:30"
$ws.Range("O9").Value = "12.7 °C"
$ws.Range("E10").Value = "2026-02-26 18:18:31"
$ws.Range("E11").Value = "2026-02-26 18:18:32"
$ws.Range("E12").Value = "2026-02-26 18:18:33"
$ws.Range("E13").Value = "2026-02-26 18:18:34"
$ws.Range("J13").Value = "1028.3 hPa"
$ws.Range("L13").Value = "27.7 km/h - 133º 17:45 TU"
$ws.Range("O13").Value = "7.1 °C"
$ws.Range("E14").Value = "2026-02-26 18:18:35"
$ws.Range("E15").Value = "2026-02-26 18:18:36"
$ws.Range("E16").Value = "2026-02-26 18:18:39"
$ws.Range("H16").Value = "'42%"
$ws.Range("E17").Value = "2026-02-26 18:18:41"
$ws.Range("G17").Value = "2 cm"
$ws.Range("E18").Value = "2026-02-26 18:18:44"
$ws.Range("E19").Value = "2026-02-26 18:18:46"
$ws.Range("H19").Value = "'44%"
$ws.Range("O19").Value = "11.6 °C"
$ws.Range("E20").Value = "2026-02-26 18:18:47"
$ws.Range("E21").Value = "2026-02-26 18:18:49"
$ws.Range("O21").Value = "9.9 °C"
$ws.Range("E22").Value = "2026-02-26 18:18:51"
$ws.Range("E23").Value = "2026-02-26 18:18:54"
$ws.Range("H23").Value = "'37%"
$ws.Range("E24").Value = "2026-02-26 18:18:56"
$ws.Range("J24").Value = "1026.8 hPa"
$ws.Range("O24").Value = "10.7 °C"
$ws.Range("E25").Value = "2026-02-26 18:18:59"
$ws.Range("H25").Value = "'35%"
$ws.Range("E26").Value = "2026-02-26 18:19:01"
$ws.Range("H26").Value = "'38%"
$ws.Range("J26").Value = "1024.1 hPa"
$ws.Range("O26").Value = "11.3 °C"
$ws.Range("E27").Value = "2026-02-26 18:19:03"
$ws.Range("E28").Value = "2026-02-26 18:19:06"
$ws.Range("H28").Value = "'77%"
$ws.Range("E29").Value = "2026-02-26 18:19:08"
$ws.Range("E30").Value = "2026-02-26 18:19:11"
$ws.Range("O30").Value = "12.5 °C"
$ws.Range("E31").Value = "2026-02-26 18:19:13"
$ws.Range("H31").Value = "'83%"
$ws.Range("E32").Value = "2026-02-26 18:19:15"
$ws.Range("H32").Value = "'61%"
$ws.Range("O32").Value = "8.6 °C"
$ws.Range("E33").Value = "2026-02-26 18:19:18"
$ws.Range("H33").Value = "'54%"
$ws.Range("J33").Value = "1026.7 hPa"
$ws.Range("O33").Value = "8.6 °C"
$ws.Range("E34").Value = "2026-02-26 18:19:20"
$ws.Range("H34").Value = "'44%"
$ws.Range("O34").Value = "5.2 °C"
$ws.Range("E35").Value = "2026-02-26 18:19:22"
$ws.Range("H35").Value = "'40%"
$ws.Range("E36").Value = "2026-02-26 18:19:25"
$ws.Range("E37").Value = "2026-02-26 18:19:27"
$ws.Range("H37").Value = "'71%"
$ws.Range("E38").Value = "2026-02-26 18:19:29"
$ws.Range("E39").Value = "2026-02-26 18:19:32"
$ws.Range("H39").Value = "'41%"
$ws.Range("N39").Value = "0.7 °C 17:33 TU"
$ws.Range("E40").Value = "2026-02-26 18:19:34"
$ws.Range("J40").Value = "1027.3 hPa"
$ws.Range("O40").Value = "9.8 °C"
$ws.Range("E41").Value = "2026-02-26 18:19:37"
$ws.Range("E42").Value = "2026-02-26 18:19:39"
$ws.Range("E43").Value = "2026-02-26 18:19:41"
$ws.Range("L43").Value = "19.1 km/h - 182º 17:44 TU"
$ws.Range("O43").Value = "9.4 °C"
$ws.Range("E44").Value = "2026-02-26 18:19:44"
$ws.Range("H44").Value = "'51%"
$ws.Range("E45").Value = "2026-02-26 18:19:46"
$ws.Range("E46").Value = "2026-02-26 18:19:49"
$ws.Range("H46").Value = "'82%"
$ws.Range("O46").Value = "11.2 °C"
